$d = $word.ActiveDocument

# Locate the end of " exceptionhandling :<" so we can split the paragraph
# right after it (this text stays in the first paragraph).
$findRange = $d.Content
$found = $findRange.Find.Execute(" exceptionhandling :<", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPos = $findRange.End

# Step 1: Split the paragraph right after " exceptionhandling :<" into two
# paragraphs. This creates a new, currently empty, paragraph right after it.
$splitPoint = $d.Range($splitPos, $splitPos)
$splitPoint.InsertParagraphAfter()

# Step 2: Insert the new sentence text into the new (second) paragraph.
$newText = "Hver anden gang man vil slette et billede i en session får vi en underlig exception i Storage.DeletePicture. Den er fanget af en try-Catch blok nu, men ikke rettet, da vi ( kewin ) ikke kunne løse det."
$newParaStart = $d.Range($splitPos + 1, $splitPos + 1)
$newParaStart.InsertAfter($newText)

# Step 3: Move the "_GoBack" bookmark (originally located between
# " omfattende" and " exceptionhandling :<") to the very end of the newly
# inserted sentence, i.e. the end of the new second paragraph.
$endOfText = $splitPos + 1 + $newText.Length

# Placing a collapsed bookmark range exactly at a paragraph's end position
# lands incorrectly in this runtime, so as a workaround we insert a temporary
# placeholder character there, anchor the bookmark around it, then delete the
# placeholder again (the bookmark stays correctly anchored at that spot).
$placeholderRange = $d.Range($endOfText, $endOfText)
$placeholderRange.InsertAfter("@")

$bmRange = $d.Range($endOfText, $endOfText + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

$delRange = $d.Range($endOfText, $endOfText + 1)
$delRange.Delete()
